$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$currencyFmt = '"$"#,##0.00;[Red]\-"$"#,##0.00'

# --- Row 3: Atmega32m1 -> Atmega324 (price/source unchanged) ---
$ws.Range("A3").Value = 'Atmega324'

# --- Row 5: Stepper - price 4.75 -> 4.63, source becomes a hyperlink ---
$ws.Range("B5").Value = 4.63
$ws.Range("C5").Value = 'http://nxtmarket.info/item/536308083756'
$ws.Hyperlinks.Add($ws.Range("C5"), 'http://nxtmarket.info/item/536308083756')

# --- Row 4: Foam Tiles - price 6 -> 4, note text updated ---
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = '(Only using a third $12 originally)'

# --- Row 6: MOSFET now has a price ---
$ws.Range("B6").Value = 0.3

# --- Rows 9-12: fill in previously-blank prices ---
$ws.Range("B9").Value = 0.6
$ws.Range("B10").Value = 0.13
$ws.Range("B11").Value = 1.1599999999999999
$ws.Range("B12").Value = 1.9

# --- Remove the "8 Pin Sockets" row (old row 13) and move Controller Housing up ---
$ws.Range("A13").Value = 'Controller Housing'
$ws.Range("B13").Value = 6.25
$ws.Range("C13").Value = 'Might get away with at $3 one, unless we want to start looking for budget options'
$ws.Range("A14:D14").Clear()

# --- Uniform currency formatting (2 decimals) across all price cells ---
$ws.Range("B2:B13").NumberFormat = $currencyFmt
$ws.Range("D2").NumberFormat = $currencyFmt

# --- Stray formatted-but-empty cell at D19 (left over from formatting drag) ---
$ws.Range("D19").NumberFormat = $currencyFmt

# --- Recompute / re-assert the total formula text ---
$ws.Range("D2").Formula = "=SUM(B2:B17)"

# --- Column B got wider ---
$ws.Columns("B").ColumnWidth = 13

# --- Selection cursor moved ---
$ws.Range("E24").Select() | Out-Null
